# Update task list: refresh statuses/assignee/priority for a few rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "3. Rename: ..." row -> Status: Open -> In Progress, Assignee -> Arthur
$ws.Range("C4").Value = "In Progress"
$ws.Range("D4").Value = "Arthur"

# "4. Include "Link Generator" as a feature" row -> Priority: Middle -> Hight
$ws.Range("B5").Value = "Hight"

# "5. Search algorithm optimization" row -> Status: In Progress -> Done (shown in green)
$ws.Range("C6").Value = "Done"
$ws.Range("C6").Font.Color = 5287936

# Keep the previously selected cell in sync with where the edits were made.
$ws.Range("E7").Select()
